# Insert a new data row at row 31 (pushing existing rows 31-38 down to 32-39)
# and populate it with a new weekly price record for "Haba" at
# "Vega Monumental Concepción".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("31:31").Insert()

$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = 44855
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 100112026
$ws.Range("G31").Value = "Haba"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 7000
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = 7500
$ws.Range("N31").Value = "$/saco 25 kilos"
$ws.Range("O31").Value = "Región Metropolitana"
$ws.Range("P31").Value = 300
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"
